$d = $word.ActiveDocument

# Simple whole-document text replacement (used where the run already
# carries explicit direct formatting, e.g. bold/italic, so the engine
# does not coalesce it with a neighbouring empty run).
function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# Some bullet paragraphs look like <w:r/><w:r><w:t>OLD</w:t></w:r> - an
# empty run immediately followed by the plain (no rPr) text run. A plain
# Find/Replace on the text run's content causes the engine to coalesce
# that leading empty run into the edited run (both have "no formatting"),
# which would silently drop the <w:r/> element. To keep the document
# shape intact, redo the replacement and then re-insert a fresh empty
# run immediately in front of the freshly edited run via InsertXML
# (which does not trigger the coalescing pass).
function Replace-Text-KeepEmptyRun($old, $new) {
    $rng = $d.Content
    $result = $rng.Find.Execute($old, $true, $false, $false, $false, $false, `
                                 $true, 1, $false, $new, 2)
    if ($result) {
        $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
        $insertPoint = $d.Range($rng.Start, $rng.Start)
        $insertPoint.InsertXML($xml)
    }
}

# Title / heading (appears twice: Heading1 at top and bold run near the
# bottom - both already carry run-level formatting info or no leading
# empty run, so the plain replace is safe for both occurrences).
Replace-Text "Play Mayan Temple Advance for Free - Review & Free Play" "Play Mayan Temple Advance - Free Slot Game"

# "What we like" bullet points
Replace-Text-KeepEmptyRun "Mayan theme and symbols are immersive and authentic" "Stunning Mayan theme with authentic symbols"
Replace-Text-KeepEmptyRun "Temple Bonus Game offers big payout opportunities" "Exciting Temple Bonus Game for big payouts"
Replace-Text-KeepEmptyRun "RTP of 95.09% indicates frequent winnings" "High RTP of 95.09% for frequent winnings"
Replace-Text-KeepEmptyRun "Gameplay similar to other popular games with a Maya theme" "Similar games available for those who enjoy the Maya theme"

# "What we don't like" bullet points
Replace-Text-KeepEmptyRun "Limited number of paylines may not appeal to some players" "Limited variety of symbols on the reels"
Replace-Text-KeepEmptyRun "Graphics are not as visually stunning as some other slot games" "No progressive jackpot feature"

# Meta description (italic run at the end - already has explicit rPr)
Replace-Text "Immerse yourself in ancient Maya civilization. Review of Mayan Temple Advance slot game with RTP, bonus features, and similar games. Play now for free." "Discover the world of the Mayans and play Mayan Temple Advance for free. Win big with the Temple Bonus Game."
